$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 39999.2
$ws.Range("J87").Value = 39999.2
$ws.Range("L87").Value = 39999.2
$ws.Range("N87").Value = -42495.2
$ws.Range("H90").Value = 39999.2
$ws.Range("J90").Value = 39999.2
$ws.Range("L90").Value = 119997.6
$ws.Range("N90").Value = -132477.6
$ws.Range("H101").Value = 400
$ws.Range("I101").Value = 400
$ws.Range("K101").Value = 1200
$ws.Range("M101").Value = 422
$ws.Range("H107").Value = 1107.375
$ws.Range("J107").Value = 2719.6
$ws.Range("L107").Value = 2719.6
$ws.Range("N107").Value = -6559.6
$ws.Range("H132").Value = 3833.4722
$ws.Range("I132").Value = 3438.7188
$ws.Range("K132").Value = 10316.1564
$ws.Range("M132").Value = -7786.1564
$ws.Range("H138").Value = 4043.348
$ws.Range("I138").Value = 2719.4
$ws.Range("J138").Value = 4411.1113
$ws.Range("K138").Value = 8158.200000000001
$ws.Range("L138").Value = 13233.3339
$ws.Range("M138").Value = -3018.200000000001
$ws.Range("N138").Value = -23513.3339
$ws.Range("H141").Value = 4389.684
$ws.Range("I141").Value = 4389.684
$ws.Range("K141").Value = 13169.052
$ws.Range("M141").Value = -7989.052

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4577.4
$ws.Range("I61").Value = 2999.5
$ws.Range("K61").Value = 2999.5
$ws.Range("M61").Value = -2787.5
$ws.Range("H74").Value = 2055.2778
$ws.Range("I74").Value = 2055.2778
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2055.2778
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = -1181.2778
$ws.Range("H77").Value = 2055.2778
$ws.Range("I77").Value = 2055.2778
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10276.389
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = -5908.388999999999
$ws.Range("H110").Value = 4500
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 4500
$ws.Range("K110").Value = 0
$ws.Range("M110").Value = 4500
$ws.Range("N110").Value = -8590
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("H136").Value = 4577.4
$ws.Range("I136").Value = 2999.5
$ws.Range("K136").Value = 8998.5
$ws.Range("M136").Value = -6448.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("H105").Value = 2308.5386
$ws.Range("I105").Value = 780.1111
$ws.Range("K105").Value = 780.1111
$ws.Range("M105").Value = 966.8889
$ws.Range("H107").Value = 2121
$ws.Range("J107").Value = 4750
$ws.Range("L107").Value = 4750
$ws.Range("N107").Value = -8590
$ws.Range("H134").Value = 14333
$ws.Range("I134").Value = 14333
$ws.Range("K134").Value = 42999
$ws.Range("M134").Value = -40464

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5813.476
$ws.Range("I58").Value = 4792.6313
$ws.Range("J58").Value = 15511.5
$ws.Range("K58").Value = 4792.6313
$ws.Range("L58").Value = 15511.5
$ws.Range("M58").Value = -4589.6313
$ws.Range("N58").Value = -15917.5
$ws.Range("H96").Value = 39000
$ws.Range("J96").Value = 39000
$ws.Range("L96").Value = 39000
$ws.Range("N96").Value = -44492
$ws.Range("H132").Value = 1849.25
$ws.Range("J132").Value = 1999.5
$ws.Range("L132").Value = 5998.5
$ws.Range("N132").Value = -11058.5
$ws.Range("H136").Value = 5813.476
$ws.Range("I136").Value = 4792.6313
$ws.Range("J136").Value = 15511.5
$ws.Range("K136").Value = 14377.8939
$ws.Range("L136").Value = 46534.5
$ws.Range("M136").Value = -11827.8939
$ws.Range("N136").Value = -51634.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1420
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 9000
$ws.Range("N11").Value = -9280
$ws.Range("H60").Value = 262.5
$ws.Range("I60").Value = 262.5
$ws.Range("K60").Value = 787.5
$ws.Range("M60").Value = -536.5
$ws.Range("H86").Value = 850
$ws.Range("H89").Value = 850
$ws.Range("H113").Value = 1699.25
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 2997
$ws.Range("M113").Value = -827
$ws.Range("H121").Value = 1826.2307
$ws.Range("I121").Value = 508.6
$ws.Range("J121").Value = 2649.75
$ws.Range("K121").Value = 1525.8
$ws.Range("L121").Value = 7949.25
$ws.Range("M121").Value = -215.8000000000002
$ws.Range("N121").Value = -10569.25
$ws.Range("H129").Value = 3956.25
$ws.Range("J129").Value = 3956.25
$ws.Range("L129").Value = 11868.75
$ws.Range("N129").Value = -21868.75
$ws.Range("H131").Value = 16572.133
$ws.Range("J131").Value = 3611.3635
$ws.Range("L131").Value = 10834.0905
$ws.Range("N131").Value = -20914.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5937.1816
$ws.Range("I70").Value = 5559.4
$ws.Range("J70").Value = 6252
$ws.Range("K70").Value = 5559.4
$ws.Range("L70").Value = 6252
$ws.Range("M70").Value = -5289.4
$ws.Range("N70").Value = -6792
$ws.Range("H73").Value = 5937.1816
$ws.Range("I73").Value = 5559.4
$ws.Range("J73").Value = 6252
$ws.Range("K73").Value = 5559.4
$ws.Range("L73").Value = 6252
$ws.Range("M73").Value = -4623.4
$ws.Range("N73").Value = -8124
$ws.Range("H113").Value = 1133
$ws.Range("I113").Value = 1133
$ws.Range("K113").Value = 1133
$ws.Range("M113").Value = 1037
$ws.Range("H132").Value = 4500
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -18560

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("H22").Value = 1812.375
$ws.Range("H27").Value = 1812.375
$ws.Range("H46").Value = 4374.375
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3312
$ws.Range("H68").Value = 3999
$ws.Range("I68").Value = 3999
$ws.Range("K68").Value = 3999
$ws.Range("M68").Value = -3250
$ws.Range("H71").Value = 3999
$ws.Range("I71").Value = 3999
$ws.Range("K71").Value = 19995
$ws.Range("M71").Value = -16251
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H100").Value = 6212.7144
$ws.Range("I100").Value = 6415.8335
$ws.Range("J100").Value = 4994
$ws.Range("K100").Value = 6415.8335
$ws.Range("L100").Value = 4994
$ws.Range("M100").Value = -5874.8335
$ws.Range("N100").Value = -6076
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 12231.523
$ws.Range("I132").Value = 14704.786
$ws.Range("J132").Value = 7285
$ws.Range("K132").Value = 44114.358
$ws.Range("L132").Value = 21855
$ws.Range("M132").Value = -41584.358
$ws.Range("N132").Value = -26915
$ws.Range("H136").Value = 6918.6665
$ws.Range("I136").Value = 3503.5
$ws.Range("J136").Value = 8626.25
$ws.Range("K136").Value = 10510.5
$ws.Range("L136").Value = 25878.75
$ws.Range("M136").Value = -7960.5
$ws.Range("N136").Value = -30978.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 44999
$ws.Range("J92").Value = 44999
$ws.Range("L92").Value = 44999
$ws.Range("N92").Value = -49991
$ws.Range("H101").Value = 20040.6
$ws.Range("J101").Value = 20040.6
$ws.Range("L101").Value = 20040.6
$ws.Range("N101").Value = -26530.6
$ws.Range("H132").Value = 2721.3635
$ws.Range("I132").Value = 2673.5
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 8020.5
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -5490.5
$ws.Range("N132").Value = -14660
$ws.Range("H136").Value = 2671
$ws.Range("I136").Value = 2792.375
$ws.Range("K136").Value = 8377.125
$ws.Range("M136").Value = -5827.125
